$d = $word.ActiveDocument

# Helper: rewrite a paragraph's text as a single clean run, stripping any
# stray w:proofErr (spell-check) markers that Word had scattered between
# the original runs (e.g. around "Qte" / "grab" or "cute" / "creepy").
#
# Plain Find/Replace across the paragraph merges the runs but leaves any
# proofErr marker that sits at the very first/last child position of the
# paragraph untouched. So we temporarily pad the paragraph with sentinel
# text on both sides (pushing every proofErr marker to an interior
# position), replace the padded text with itself (which merges everything,
# including the now-interior proofErr markers, into one run), and finally
# strip the sentinel padding back off.
function Clear-ParaProofErr($paraIndex) {
    $p = $d.Paragraphs($paraIndex)
    $start = $p.Range.Start
    $endExclMark = $p.Range.End - 1
    $original = $d.Range($start, $endExclMark).Text

    $padStart = [string][char]0xE000 + "PAD" + [string][char]0xE000
    $padEnd = [string][char]0xE000 + "DAP" + [string][char]0xE000

    $d.Range($start, $start).InsertBefore($padStart) | Out-Null
    $p2 = $d.Paragraphs($paraIndex)
    $endExclMark2 = $p2.Range.End - 1
    $d.Range($endExclMark2, $endExclMark2).InsertAfter($padEnd) | Out-Null

    $paddedText = $padStart + $original + $padEnd
    $d.Content.Find.Execute($paddedText, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $paddedText, 2) | Out-Null
    $d.Content.Find.Execute($paddedText, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $original, 2) | Out-Null
}

# 1) "Qte sorcière grab" paragraph (2nd paragraph in the doc).
Clear-ParaProofErr 2

# 2) "Ils voient deux chemins devant eux, un cute un creepy. CHOIX" paragraph.
$targetIdx = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -match "deux chemins devant eux") {
        $targetIdx = $i
    }
}
Clear-ParaProofErr $targetIdx

# 3) Add the two new paragraphs + "EWAN" run right before the trailing
#    bookmark-only paragraph.
$last = $d.Paragraphs($d.Paragraphs.Count)
$last.Range.InsertParagraphBefore()
$p1 = $d.Paragraphs($d.Paragraphs.Count - 1)
$p1.Range.Text = "Ils voient une maison faite de bonbons et de pain d’épices, ils sont attirés par cette maison et s’avancent."

$last2 = $d.Paragraphs($d.Paragraphs.Count)
$last2.Range.InsertParagraphBefore()
$p2 = $d.Paragraphs($d.Paragraphs.Count - 1)
$p2.Range.Text = "Porte interactive pour toquer, temps de latence avant que la sorcière ouvre"

$bookmarkPara = $d.Paragraphs($d.Paragraphs.Count)
$ins = $d.Range($bookmarkPara.Range.Start, $bookmarkPara.Range.Start)
$ins.InsertBefore("EWAN")
